$wb = $excel.ActiveWorkbook

# Update "想去人数" (column F) values on both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F10").Value = 4985
    $ws.Range("F11").Value = 4686
    $ws.Range("F15").Value = 40
}
